# Apply the betexplorer Eliteserien 2023 update:
#  - Rows 197-200 get re-ordered/re-synced match data (same 4 matches, new arrangement)
#  - A brand-new match row (217) is appended at the bottom of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rewrite the data for rows 197-200 (F..V). Columns A-D stay as they are
#    (same index/country/tournament/season) and column E (kickoff date/time)
#    is identical for all four of these matches already, so only F..V need
#    new values.
# ---------------------------------------------------------------------------

# Row 197 -> Stromsgodset 3-1 Odd
$ws.Range("F197").Value = "Stromsgodset"
$ws.Range("G197").Value = 3
$ws.Range("H197").Value = "Odd"
$ws.Range("I197").Value = 1
$ws.Range("J197").Value = 2.15
$ws.Range("K197").Value = "09/10/2023 14:42"
$ws.Range("L197").Value = 1.97
$ws.Range("M197").Value = "22/10/2023 16:53"
$ws.Range("N197").Value = 3.63
$ws.Range("O197").Value = "09/10/2023 14:42"
$ws.Range("P197").Value = 3.56
$ws.Range("Q197").Value = "22/10/2023 16:54"
$ws.Range("R197").Value = 3.42
$ws.Range("S197").Value = "09/10/2023 14:42"
$ws.Range("T197").Value = 4.14
$ws.Range("U197").Value = "22/10/2023 16:53"
$ws.Range("V197").Value = "https://www.betexplorer.com/football/norway/eliteserien/stromsgodset-odds-bk/EaUYt1Zk/"

# Row 198 -> Rosenborg 1-1 Stabaek
$ws.Range("F198").Value = "Rosenborg"
$ws.Range("G198").Value = 1
$ws.Range("H198").Value = "Stabaek"
$ws.Range("I198").Value = 1
$ws.Range("J198").Value = 1.79
$ws.Range("K198").Value = "09/10/2023 14:42"
$ws.Range("L198").Value = 1.98
$ws.Range("M198").Value = "22/10/2023 16:56"
$ws.Range("N198").Value = 4.01
$ws.Range("O198").Value = "09/10/2023 14:42"
$ws.Range("P198").Value = 3.88
$ws.Range("Q198").Value = "22/10/2023 16:53"
$ws.Range("R198").Value = 4.27
$ws.Range("S198").Value = "09/10/2023 14:42"
$ws.Range("T198").Value = 3.75
$ws.Range("U198").Value = "22/10/2023 16:56"
$ws.Range("V198").Value = "https://www.betexplorer.com/football/norway/eliteserien/rosenborg-stabaek/KxVUssKq/"

# Row 199 -> HamKam 0-3 Haugesund
$ws.Range("F199").Value = "HamKam"
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = "Haugesund"
$ws.Range("I199").Value = 3
$ws.Range("J199").Value = 2.05
$ws.Range("K199").Value = "09/10/2023 14:42"
$ws.Range("L199").Value = 2.28
$ws.Range("M199").Value = "22/10/2023 16:55"
$ws.Range("N199").Value = 3.74
$ws.Range("O199").Value = "09/10/2023 14:42"
$ws.Range("P199").Value = 3.52
$ws.Range("Q199").Value = "22/10/2023 16:55"
$ws.Range("R199").Value = 3.51
$ws.Range("S199").Value = "09/10/2023 14:42"
$ws.Range("T199").Value = 3.29
$ws.Range("U199").Value = "22/10/2023 16:55"
$ws.Range("V199").Value = "https://www.betexplorer.com/football/norway/eliteserien/ham-kam-haugesund/j7C6m3dM/"

# Row 200 -> Aalesund 3-2 Sarpsborg 08
$ws.Range("F200").Value = "Aalesund"
$ws.Range("G200").Value = 3
$ws.Range("H200").Value = "Sarpsborg 08"
$ws.Range("I200").Value = 2
$ws.Range("J200").Value = 3.9
$ws.Range("K200").Value = "09/10/2023 14:42"
$ws.Range("L200").Value = 3.58
$ws.Range("M200").Value = "22/10/2023 16:53"
$ws.Range("N200").Value = 4.25
$ws.Range("O200").Value = "09/10/2023 14:42"
$ws.Range("P200").Value = 4.21
$ws.Range("Q200").Value = "22/10/2023 16:53"
$ws.Range("R200").Value = 1.82
$ws.Range("S200").Value = "09/10/2023 14:42"
$ws.Range("T200").Value = 1.95
$ws.Range("U200").Value = "22/10/2023 16:53"
$ws.Range("V200").Value = "https://www.betexplorer.com/football/norway/eliteserien/aalesund-sarpsborg-08/4jAfj5C3/"

# ---------------------------------------------------------------------------
# 2) Append the new match row (217) at the bottom, reusing the formatting
#    of the previous last row (216) so the bold/bordered index column and
#    the date-time number format on column E carry over correctly.
# ---------------------------------------------------------------------------

$ws.Range("A216:V216").Copy($ws.Range("A217:V217"))

$ws.Range("A217").Value = 216
$ws.Range("B217").Value = "norway"
$ws.Range("C217").Value = "eliteserien"

# "2023" looks numeric, so a plain assignment would silently become the
# number 2023 instead of text (like every other row's "temporada" column).
# Force text by switching to a Text format, assigning, then pasting the
# General format back over it (value - a shared string - is left alone).
$ws.Range("D217").NumberFormat = "@"
$ws.Range("D217").Value = "2023"
$ws.Range("D216").Copy()
$ws.Range("D217").PasteSpecial(-4122)

$ws.Range("E217").Value = 45236.791666666664
$ws.Range("F217").Value = "Aalesund"
$ws.Range("G217").Value = 0
$ws.Range("H217").Value = "Sandefjord"
$ws.Range("I217").Value = 3
$ws.Range("J217").Value = 2.19
$ws.Range("K217").Value = "30/10/2023 19:13"
$ws.Range("L217").Value = 2.07
$ws.Range("M217").Value = "06/11/2023 18:59"
$ws.Range("N217").Value = 3.77
$ws.Range("O217").Value = "30/10/2023 19:13"
$ws.Range("P217").Value = 4.01
$ws.Range("Q217").Value = "06/11/2023 18:59"
$ws.Range("R217").Value = 3.15
$ws.Range("S217").Value = "30/10/2023 19:13"
$ws.Range("T217").Value = 3.38
$ws.Range("U217").Value = "06/11/2023 18:59"
$ws.Range("V217").Value = "https://www.betexplorer.com/football/norway/eliteserien/aalesund-sandefjord/OUmHgvJe/"

Write-Output "done"
